# Sigi-11: Remove coordinates to json logic from C# script
#
# The spreadsheet used to carry two separate rows for the geographic
# coordinates of each seal (row 106 = LATITUDE, row 107 = LONGITUDE), one
# value per seal column (B = SigiDoc 13, C = SigiDoc 11). Since the C#
# importer no longer turns these into a lat/long JSON pair, the two rows
# are collapsed into a single "COORDINATES" row containing the combined
# "lat, long" string per seal, and the now-empty LONGITUDE row is cleared
# out (only the styled, value-less A cell is left behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the combined "lat, long" strings first (columns B/C), then the new
# row label in column A, so that new shared-string entries are appended in
# that same order.
$ws.Range("B106").Value = "42.136097, 24.742168"
$ws.Range("C106").Value = "42.698334, 23.319941"
$ws.Range("A106").Value = "COORDINATES"

# Row 107 (the old LONGITUDE row) loses its data entirely. B107/C107 are
# fully cleared (value + formatting) so no cell element remains for them;
# A107 keeps its row-label styling but its text is removed.
$ws.Range("B107:C107").Clear()
$ws.Range("A107").ClearContents()

# Reflect where the user ended up after the edit: scrolled down near the
# bottom of the sheet with the now-blank A107 selected.
$ws.Range("A107").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
